# "A lot of spelling later" - retitle the M6 poster-event poster.
#
# Title shape ("Title 1", the first shape on slide 1) gets its text
# reworded and its run's language tag updated from en-US to en-GB.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$titleShape = $s.Shapes.Title
$titleRange = $titleShape.TextFrame.TextRange

$titleRange.Text = "Fuzz Testing of Constraint Programming"
$titleRange.Runs(1).LanguageID = "en-GB"
